# "Feito 1 e 2" — mark exercises 1 and 2 (rows 2 and 3) as done by
# checking the boolean "done" column (C), matching the existing
# checkboxes already present from row 4 downward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $true

# Leave the cursor/selection on A4, which is where it ended up after
# checking those boxes.
$ws.Range("A4").Select()
